$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data row (row 2): name, email (hyperlink), repo link (hyperlink)
$ws.Range("A2").Value = "محمد حسين عبدالحافظ محمد "

$ws.Range("B2").Value = "mhbadawi14@gmail.com"
$ws.Hyperlinks.Add($ws.Range("B2"), "mailto:mhbadawi14@gmail.com") | Out-Null

$ws.Range("C2").Value = "https://github.com/mhbadawi14/oos"
$ws.Hyperlinks.Add($ws.Range("C2"), "https://github.com/mhbadawi14/oos") | Out-Null

# Match the saved selection in the target workbook
$ws.Range("E6").Select() | Out-Null
